# Applies the "Deploying to gh-pages" content refresh to the FHIR
# StructureDefinition-rating-area workbook:
#   - Metadata sheet: bump Version, refresh publish Date, add a Publisher
#     name, and replace the (duplicated) "Contact" property row with a new
#     "Jurisdiction" property row.
#   - Elements sheet: the generated "Short"/"Definition" text for the root
#     Extension row now reflects the resource's own title/description
#     instead of the generic boilerplate.

$wb = $excel.ActiveWorkbook

# ---- Metadata sheet -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# Version 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date refreshed to the new publish timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# Remove one of the duplicate "Contact" rows (row 10); the remaining
# duplicate (now shifted up into row 10) is turned into the new
# "Jurisdiction" property.
$meta.Range("A10").EntireRow.Delete()
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---- Elements sheet --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# Root Extension row: Short/Definition now driven by the resource's own
# Title/Description rather than the generic "Extension"/"An Extension".
$elements.Range("K2").Value = "Rating Area"
$elements.Range("L2").Value = "Code for the geographic insurance rating area of the associated plan"
